$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The decision-table rule row 11 ("R40") is being re-keyed to "1": the
# cell keeps its existing formatting (border, etc.) but must hold the
# literal text "1" rather than the numeric value 1, so force the cell
# to Text format before writing it (otherwise "1" is auto-coerced to a
# number by Excel).
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
